# Übertragung der Parameter in ein Excel Sheet
# Adds new survey submission rows received 2024-12-09 (evening) to the
# various input tables of the workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# contact_data: two new (mostly empty) contact submissions
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("contact_data")
$ws.Range("A20").Value = 1
$ws.Range("B20").Value = "2024-12-09 19:17:54"
$ws.Range("E20").Value = "SKZ"

$ws.Range("A21").Value = 1
$ws.Range("B21").Value = "2024-12-09 19:38:23"
$ws.Range("E21").Value = "SKZ"

# ---------------------------------------------------------------------
# company_data: two new company submissions (same address, repeated)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("company_data")

$ws.Range("A21").Value = 1
$ws.Range("B21").Value = "2024-12-09 19:18:15"
$ws.Range("C21").Value = "SKZ"
$ws.Range("D21").Value = "Friedrich-Bergius-Ring 22"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "97076"
$ws.Range("F21").Value = "Würzburg"
$ws.Range("G21").Value = "Bayern"
$ws.Range("H21").Value = "Deutschland"
$ws.Range("I21").Value = $false
$ws.Range("J21").Value = $false
$ws.Range("L21").Value = 49.80282025
$ws.Range("M21").Value = 10.00010726291456

$ws.Range("A22").Value = 1
$ws.Range("B22").Value = "2024-12-09 19:38:44"
$ws.Range("C22").Value = "SKZ"
$ws.Range("D22").Value = "Friedrich-Bergius-Ring 22"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "97076"
$ws.Range("F22").Value = "Würzburg"
$ws.Range("G22").Value = "Bayern"
$ws.Range("H22").Value = "Deutschland"
$ws.Range("I22").Value = $false
$ws.Range("J22").Value = $false
$ws.Range("L22").Value = 49.80282025
$ws.Range("M22").Value = 10.00010726291456

# ---------------------------------------------------------------------
# product_fractions: four new fraction submissions
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("product_fractions")

$ws.Range("A22").Value = 1
$ws.Range("B22").Value = "2024-12-09 19:18:43"
$ws.Range("C22").Value = "['PE-LD', 'PP', 'PVC-U', 'SAN']"
$ws.Range("D22").Value = "['', '', '', '']"
$ws.Range("E22").Value = "[50.0, 30.0, 15.0, 5.0]"

$ws.Range("A23").Value = 1
$ws.Range("B23").Value = "2024-12-09 19:39:03"
$ws.Range("C23").Value = "['PE-LD', 'PP', 'PVC-U', 'POM']"
$ws.Range("D23").Value = "['', '', '', '']"
$ws.Range("E23").Value = "[50.0, 30.0, 10.0, 10.0]"

$ws.Range("A24").Value = 1
$ws.Range("B24").Value = "2024-12-09 19:40:30"
$ws.Range("C24").Value = "['PE-LD', 'PP', 'PVC-U', 'Aluminium']"
$ws.Range("D24").Value = "['', '', '', '']"
$ws.Range("E24").Value = "[50.0, 30.0, 10.0, 10.0]"

$ws.Range("A25").Value = 1
$ws.Range("B25").Value = "2024-12-09 19:50:28"
$ws.Range("C25").Value = "['PP', 'PS']"
$ws.Range("D25").Value = "['', '']"
$ws.Range("E25").Value = "[50.0, 50.0]"

# ---------------------------------------------------------------------
# product_origin: one new submission
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("product_origin")

$ws.Range("A11").Value = 1
$ws.Range("B11").Value = "2024-12-09 19:18:47"
$ws.Range("C11").Value = "Post-Industrial (PI)"

# ---------------------------------------------------------------------
# product_quality: three new submissions
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("product_quality")

$ws.Range("A20").Value = 1
$ws.Range("B20").Value = "2024-12-09 19:19:00"
$ws.Range("C20").Value = "Ja"
$ws.Range("D20").Value = "braun"
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = "mittel"
$ws.Range("H20").Value = "[[], [], [], []]"
$ws.Range("I20").Value = "[[], [], [], []]"

$ws.Range("A21").Value = 1
$ws.Range("B21").Value = "2024-12-09 19:39:15"
$ws.Range("C21").Value = "Ja"
$ws.Range("D21").Value = "grün"
$ws.Range("E21").Value = 99.98999999999999
$ws.Range("F21").Value = "gering"
$ws.Range("H21").Value = "[]"
$ws.Range("I21").Value = "[]"

$ws.Range("A22").Value = 1
$ws.Range("B22").Value = "2024-12-09 19:50:40"
$ws.Range("C22").Value = "Ja"
$ws.Range("E22").Value = 100
$ws.Range("F22").Value = "gering"
$ws.Range("H22").Value = "[]"
$ws.Range("I22").Value = "[]"

# ---------------------------------------------------------------------
# product_amount: three new submissions
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("product_amount")

$ws.Range("A20").Value = 1
$ws.Range("B20").Value = "2024-12-09 19:19:10"
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 10
$ws.Range("E20").Value = "Quartal"

$ws.Range("A21").Value = 1
$ws.Range("B21").Value = "2024-12-09 19:39:24"
$ws.Range("C21").Value = 5
$ws.Range("D21").Value = 10
$ws.Range("E21").Value = "Jahr"

$ws.Range("A22").Value = 1
$ws.Range("B22").Value = "2024-12-09 19:50:49"
$ws.Range("C22").Value = 5
$ws.Range("D22").Value = 5
$ws.Range("E22").Value = "Woche"
